$wb = $excel.ActiveWorkbook

# --- Sheet "M" ---
$wsM = $wb.Worksheets.Item("M")
$wsM.Range("E4").Value = 1.75
$wsM.Activate()
$wsM.Range("E7").Select()

# --- Sheet "Y" ---
$wsY = $wb.Worksheets.Item("Y")
$wsY.Range("G2").Value = 100
$wsY.Range("C6").Value = 0.5
$wsY.Range("C7").Value = 0.5
$wsY.Activate()
$wsY.Range("H10").Select()
